$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474993058530/1743637847690/PaZ2GKc6aN-camera-video-segment-driver-1743637850190.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSK7PGOQ3S%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T190200Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIDJmYJCm4w8ze5QZx7jAwT9fHpIVGsbWKlTo24nY9sUQAiBZbxgURLnhB%2Bjv4M26qZ8Ndj6fKVO7WNsuMVVKbUJajCrmAwjx%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMEagwDbo9An1MuK%2FMKroDHmEwCyibiRECyB1N6vtSGDYSvCmyQykvhOYJPzZ3mWgp7ei9gpVC9i5JR%2BZ%2BTVcY7pymQC8%2BOJlMcVEVUiNk%2FbARctwF8ofBDXFlWGeD9h1BYnUrHq3bC6c3uuot6fpEHwgT1TEklpQ8oRycI2OZjYAYafR8ETOQeDDH%2B8qRYNic1UKsy9FG26aupV8o4ceCZr%2FVM7YgrYju3anK8%2FrtmTIXrrZAcJZo0yxfdNJ7xY4ANriYMGgTNOm7E27nowU2uZmLjJhoICE5jDzU8MCmrbrmAGYVc1aY3eYFPS7s3UZiTa3s1CaW5e4sQkFDcXpQ9tbjya9azrBOo67p7aDC6bOvyh48118gL%2BHAGQ7fF74oTGKaR6ZbHr%2FRJYs0llYlri9lQyWnBQtk69mf3Bq5IKH2hQvT%2FVSccwliOp7fbnv8LIcKLVEb7EzJCzOfaPEXfuDXY0m7x89unVQgfjvQJHjLMfRLkmWgkHscaGLvtfq9tb2sUxfCO1uqB0JoSL3hmKUp0U9UotBxKFKpybJpHmu94vUEryEkE%2FDyxEW3qu4FM0lNPqqZfX4jXC7nmRJ2xsNVjFlXL%2BY0TDCI2Lq%2FBjqmAcDttX9dGCITq00yC%2F0yMb4QVno%2BjMB7LN9%2BY0qG5%2BbicmsypvNA6vf35OIRdaLIZwogR5HxV5nRiqE%2FuxMDAnrVIhByI4tvZ8lm%2FG1D8eZ0DbXO1AOJFVKRTKRmLPzX9tksDvgFsbD%2FujGT9Kia0kfxONir%2BEIq%2FPO700KxSKBzvwJr65XmJpw1Fm15soi3I82DegFgmjJTFfNFkBt2MhyVZHJZJxA%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2003%3A02%3A00%20GMT&X-Amz-Signature=c8ee33bf7ba5b114baa7016fb0d951e71d7a6b485657f06333a537e247b31043'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474993058530/1743634721601/5dyCDJLTC3-camera-video-segment-driver-1743634724101.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSK7PGOQ3S%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T190200Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIDJmYJCm4w8ze5QZx7jAwT9fHpIVGsbWKlTo24nY9sUQAiBZbxgURLnhB%2Bjv4M26qZ8Ndj6fKVO7WNsuMVVKbUJajCrmAwjx%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMEagwDbo9An1MuK%2FMKroDHmEwCyibiRECyB1N6vtSGDYSvCmyQykvhOYJPzZ3mWgp7ei9gpVC9i5JR%2BZ%2BTVcY7pymQC8%2BOJlMcVEVUiNk%2FbARctwF8ofBDXFlWGeD9h1BYnUrHq3bC6c3uuot6fpEHwgT1TEklpQ8oRycI2OZjYAYafR8ETOQeDDH%2B8qRYNic1UKsy9FG26aupV8o4ceCZr%2FVM7YgrYju3anK8%2FrtmTIXrrZAcJZo0yxfdNJ7xY4ANriYMGgTNOm7E27nowU2uZmLjJhoICE5jDzU8MCmrbrmAGYVc1aY3eYFPS7s3UZiTa3s1CaW5e4sQkFDcXpQ9tbjya9azrBOo67p7aDC6bOvyh48118gL%2BHAGQ7fF74oTGKaR6ZbHr%2FRJYs0llYlri9lQyWnBQtk69mf3Bq5IKH2hQvT%2FVSccwliOp7fbnv8LIcKLVEb7EzJCzOfaPEXfuDXY0m7x89unVQgfjvQJHjLMfRLkmWgkHscaGLvtfq9tb2sUxfCO1uqB0JoSL3hmKUp0U9UotBxKFKpybJpHmu94vUEryEkE%2FDyxEW3qu4FM0lNPqqZfX4jXC7nmRJ2xsNVjFlXL%2BY0TDCI2Lq%2FBjqmAcDttX9dGCITq00yC%2F0yMb4QVno%2BjMB7LN9%2BY0qG5%2BbicmsypvNA6vf35OIRdaLIZwogR5HxV5nRiqE%2FuxMDAnrVIhByI4tvZ8lm%2FG1D8eZ0DbXO1AOJFVKRTKRmLPzX9tksDvgFsbD%2FujGT9Kia0kfxONir%2BEIq%2FPO700KxSKBzvwJr65XmJpw1Fm15soi3I82DegFgmjJTFfNFkBt2MhyVZHJZJxA%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2003%3A02%3A00%20GMT&X-Amz-Signature=2c1ebf64f62b111d557f9af9b9a1ee320393e55217a537dce23c69c32e61615c'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109864/1743618603072/lsfY7Elc8s-camera-video-segment-driver-1743618608072.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSK7PGOQ3S%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T190200Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIDJmYJCm4w8ze5QZx7jAwT9fHpIVGsbWKlTo24nY9sUQAiBZbxgURLnhB%2Bjv4M26qZ8Ndj6fKVO7WNsuMVVKbUJajCrmAwjx%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMEagwDbo9An1MuK%2FMKroDHmEwCyibiRECyB1N6vtSGDYSvCmyQykvhOYJPzZ3mWgp7ei9gpVC9i5JR%2BZ%2BTVcY7pymQC8%2BOJlMcVEVUiNk%2FbARctwF8ofBDXFlWGeD9h1BYnUrHq3bC6c3uuot6fpEHwgT1TEklpQ8oRycI2OZjYAYafR8ETOQeDDH%2B8qRYNic1UKsy9FG26aupV8o4ceCZr%2FVM7YgrYju3anK8%2FrtmTIXrrZAcJZo0yxfdNJ7xY4ANriYMGgTNOm7E27nowU2uZmLjJhoICE5jDzU8MCmrbrmAGYVc1aY3eYFPS7s3UZiTa3s1CaW5e4sQkFDcXpQ9tbjya9azrBOo67p7aDC6bOvyh48118gL%2BHAGQ7fF74oTGKaR6ZbHr%2FRJYs0llYlri9lQyWnBQtk69mf3Bq5IKH2hQvT%2FVSccwliOp7fbnv8LIcKLVEb7EzJCzOfaPEXfuDXY0m7x89unVQgfjvQJHjLMfRLkmWgkHscaGLvtfq9tb2sUxfCO1uqB0JoSL3hmKUp0U9UotBxKFKpybJpHmu94vUEryEkE%2FDyxEW3qu4FM0lNPqqZfX4jXC7nmRJ2xsNVjFlXL%2BY0TDCI2Lq%2FBjqmAcDttX9dGCITq00yC%2F0yMb4QVno%2BjMB7LN9%2BY0qG5%2BbicmsypvNA6vf35OIRdaLIZwogR5HxV5nRiqE%2FuxMDAnrVIhByI4tvZ8lm%2FG1D8eZ0DbXO1AOJFVKRTKRmLPzX9tksDvgFsbD%2FujGT9Kia0kfxONir%2BEIq%2FPO700KxSKBzvwJr65XmJpw1Fm15soi3I82DegFgmjJTFfNFkBt2MhyVZHJZJxA%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2003%3A02%3A00%20GMT&X-Amz-Signature=637364a878e266536e446e2d5dda9ad5e3895b171c6d3926a54c02dd2879870c'
$ws.Range("L5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743618603072/krpw3PxSiC-camera-video-segment-1743618608072.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSK7PGOQ3S%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T190200Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIDJmYJCm4w8ze5QZx7jAwT9fHpIVGsbWKlTo24nY9sUQAiBZbxgURLnhB%2Bjv4M26qZ8Ndj6fKVO7WNsuMVVKbUJajCrmAwjx%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMEagwDbo9An1MuK%2FMKroDHmEwCyibiRECyB1N6vtSGDYSvCmyQykvhOYJPzZ3mWgp7ei9gpVC9i5JR%2BZ%2BTVcY7pymQC8%2BOJlMcVEVUiNk%2FbARctwF8ofBDXFlWGeD9h1BYnUrHq3bC6c3uuot6fpEHwgT1TEklpQ8oRycI2OZjYAYafR8ETOQeDDH%2B8qRYNic1UKsy9FG26aupV8o4ceCZr%2FVM7YgrYju3anK8%2FrtmTIXrrZAcJZo0yxfdNJ7xY4ANriYMGgTNOm7E27nowU2uZmLjJhoICE5jDzU8MCmrbrmAGYVc1aY3eYFPS7s3UZiTa3s1CaW5e4sQkFDcXpQ9tbjya9azrBOo67p7aDC6bOvyh48118gL%2BHAGQ7fF74oTGKaR6ZbHr%2FRJYs0llYlri9lQyWnBQtk69mf3Bq5IKH2hQvT%2FVSccwliOp7fbnv8LIcKLVEb7EzJCzOfaPEXfuDXY0m7x89unVQgfjvQJHjLMfRLkmWgkHscaGLvtfq9tb2sUxfCO1uqB0JoSL3hmKUp0U9UotBxKFKpybJpHmu94vUEryEkE%2FDyxEW3qu4FM0lNPqqZfX4jXC7nmRJ2xsNVjFlXL%2BY0TDCI2Lq%2FBjqmAcDttX9dGCITq00yC%2F0yMb4QVno%2BjMB7LN9%2BY0qG5%2BbicmsypvNA6vf35OIRdaLIZwogR5HxV5nRiqE%2FuxMDAnrVIhByI4tvZ8lm%2FG1D8eZ0DbXO1AOJFVKRTKRmLPzX9tksDvgFsbD%2FujGT9Kia0kfxONir%2BEIq%2FPO700KxSKBzvwJr65XmJpw1Fm15soi3I82DegFgmjJTFfNFkBt2MhyVZHJZJxA%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2003%3A02%3A00%20GMT&X-Amz-Signature=ff1f8d91eb9dc96429965f0c5b25142c8d3582ad5f6616cdf0e31f182edaf7e1'
